$d = $word.ActiveDocument

# Locate the two consecutive "Week 3" placeholder paragraphs: they are the
# two paragraphs immediately following the "Documentations have been
# revised..." paragraph, each currently containing 13 spaces.
$count = $d.Paragraphs.Count
$firstIndex = 0
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Documentations have been revised and errors have been corrected\.") {
        $firstIndex = $i + 1
        break
    }
}

# First placeholder paragraph becomes the "Week 3" summary line.
$p1 = $d.Paragraphs.Item($firstIndex)
$p1.Range.Text = "Week 3: - VGA HEX Display which has been generated from Python code has been added to the main file of the project. "

# Insert a new paragraph after it for the second bullet.
$p1 = $d.Paragraphs.Item($firstIndex)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($firstIndex + 1)
$p2.Range.Text = "              - Fixed a problem where x and y coordinates in the VGA HEX Display module were generated by the Python program from left to right in order"

# Insert another new paragraph after that for the third bullet.
$p2 = $d.Paragraphs.Item($firstIndex + 1)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($firstIndex + 2)
$p3.Range.Text = "              - Added offset feature to the HEX Display module to allow the HEX Display panel to be plotted anywhere on the screen as long as sufficient pixel space is available at that coordinate"

# The paragraph that used to be the second 13-space placeholder (now shifted
# down by two) keeps acting as a spacer, but its text shrinks from 13 to 12
# spaces.
$p4 = $d.Paragraphs.Item($firstIndex + 3)
$p4.Range.Text = "            "
